$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-01-13 Saturday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-01-14 Sunday", 2) | Out-Null
$d.Content.Find.Execute("463÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "524÷9=", 2) | Out-Null
$d.Content.Find.Execute("829÷8=", $true, $true, $false, $false, $false, $true, 1, $false, "612÷8=", 2) | Out-Null
$d.Content.Find.Execute("579÷6=", $true, $true, $false, $false, $false, $true, 1, $false, "339÷2=", 2) | Out-Null
$d.Content.Find.Execute("702÷4=", $true, $true, $false, $false, $false, $true, 1, $false, "972÷6=", 2) | Out-Null
$d.Content.Find.Execute("822÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "444÷8=", 2) | Out-Null
$d.Content.Find.Execute("503÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "878÷5=", 2) | Out-Null
$d.Content.Find.Execute("854÷9=", $true, $true, $false, $false, $false, $true, 1, $false, "364÷8=", 2) | Out-Null
$d.Content.Find.Execute("231÷3=", $true, $true, $false, $false, $false, $true, 1, $false, "662÷2=", 2) | Out-Null
$d.Content.Find.Execute("375÷9=", $true, $true, $false, $false, $false, $true, 1, $false, "529÷8=", 2) | Out-Null
$d.Content.Find.Execute("712÷3=", $true, $true, $false, $false, $false, $true, 1, $false, "227÷7=", 2) | Out-Null
$d.Content.Find.Execute("657÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "513÷9=", 2) | Out-Null
$d.Content.Find.Execute("840÷4=", $true, $true, $false, $false, $false, $true, 1, $false, "657÷3=", 2) | Out-Null
$d.Content.Find.Execute("658÷8=", $true, $true, $false, $false, $false, $true, 1, $false, "148÷5=", 2) | Out-Null
$d.Content.Find.Execute("323÷2=", $true, $true, $false, $false, $false, $true, 1, $false, "356÷5=", 2) | Out-Null
$d.Content.Find.Execute("915÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "836÷7=", 2) | Out-Null
$d.Content.Find.Execute("455÷5=", $true, $true, $false, $false, $false, $true, 1, $false, "364÷6=", 2) | Out-Null
$d.Content.Find.Execute("232÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "930÷8=", 2) | Out-Null
$d.Content.Find.Execute("470÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "521÷5=", 2) | Out-Null
$d.Content.Find.Execute("463÷9=", $true, $true, $false, $false, $false, $true, 1, $false, "519÷9=", 2) | Out-Null
$d.Content.Find.Execute("115÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "589÷4=", 2) | Out-Null
$d.Content.Find.Execute("305÷9=", $true, $true, $false, $false, $false, $true, 1, $false, "274÷3=", 2) | Out-Null
$d.Content.Find.Execute("724÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "499÷3=", 2) | Out-Null
$d.Content.Find.Execute("351÷7=", $true, $true, $false, $false, $false, $true, 1, $false, "689÷9=", 2) | Out-Null
$d.Content.Find.Execute("732÷6=", $true, $true, $false, $false, $false, $true, 1, $false, "236÷2=", 2) | Out-Null
$d.Content.Find.Execute("328÷4=", $true, $true, $false, $false, $false, $true, 1, $false, "742÷2=", 2) | Out-Null

$d.Save()
